$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-28 Sunday" "2024-01-29 Monday"

Replace-Text "469×4=" "646×5="
Replace-Text "383×2=" "117×9="
Replace-Text "754×5=" "545×6="
Replace-Text "902×8=" "671×9="
Replace-Text "817×9=" "266×2="
Replace-Text "644×3=" "386×8="
Replace-Text "161×9=" "399×9="
Replace-Text "591×3=" "706×5="
Replace-Text "541×8=" "879×3="
Replace-Text "975×4=" "903×2="
Replace-Text "133×5=" "273×7="
Replace-Text "251×9=" "397×3="
Replace-Text "915×7=" "211×9="
Replace-Text "491×8=" "678×4="
Replace-Text "193×9=" "841×7="
Replace-Text "817×4=" "868×6="
Replace-Text "893×9=" "179×7="
Replace-Text "810×7=" "820×6="
Replace-Text "185×3=" "188×3="
Replace-Text "160×7=" "358×9="
Replace-Text "617×8=" "931×7="
Replace-Text "880×4=" "719×4="
Replace-Text "482×4=" "909×3="
Replace-Text "633×9=" "199×7="
Replace-Text "423×8=" "443×5="
